$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4142.857
$ws.Range("I64").Value = 4029.4119
$ws.Range("J64").Value = 4625
$ws.Range("K64").Value = 4029.4119
$ws.Range("L64").Value = 4625
$ws.Range("M64").Value = -3781.4119
$ws.Range("N64").Value = -5121
$ws.Range("H67").Value = 4142.857
$ws.Range("I67").Value = 4029.4119
$ws.Range("J67").Value = 4625
$ws.Range("K67").Value = 4029.4119
$ws.Range("L67").Value = 4625
$ws.Range("M67").Value = -3171.4119
$ws.Range("N67").Value = -6341
$ws.Range("H138").Value = 4769.8
$ws.Range("I138").Value = 1085.375
$ws.Range("J138").Value = 8170.8076
$ws.Range("K138").Value = 3256.125
$ws.Range("L138").Value = 24512.4228
$ws.Range("M138").Value = 1883.875
$ws.Range("N138").Value = -34792.4228
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 70006
$ws.Range("I23").Value = 70006
$ws.Range("K23").Value = 70006
$ws.Range("M23").Value = -69747
$ws.Range("H32").Value = 5202.643
$ws.Range("I32").Value = 3806.6086
$ws.Range("K32").Value = 3806.6086
$ws.Range("M32").Value = -3519.6086
$ws.Range("H37").Value = 5595
$ws.Range("I37").Value = 5595
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5595
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -5322
$ws.Range("H44").Value = 21995
$ws.Range("J44").Value = 21995
$ws.Range("L44").Value = 21995
$ws.Range("N44").Value = -22971
$ws.Range("H55").Value = 16742.5
$ws.Range("J55").Value = 21990
$ws.Range("L55").Value = 21990
$ws.Range("N55").Value = -22620
$ws.Range("H64").Value = 39000
$ws.Range("J64").Value = 39000
$ws.Range("L64").Value = 39000
$ws.Range("N64").Value = -39496
$ws.Range("H67").Value = 39000
$ws.Range("J67").Value = 39000
$ws.Range("L67").Value = 39000
$ws.Range("N67").Value = -40716
$ws.Range("H80").Value = 10100
$ws.Range("I80").Value = 10100
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 10100
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -9102
$ws.Range("H83").Value = 10100
$ws.Range("I83").Value = 10100
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 30300
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -25308
$ws.Range("H102").Value = 4631280.5
$ws.Range("I102").Value = 4631280.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4631280.5
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -4629658.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23750.412
$ws.Range("I82").Value = 3652.3333
$ws.Range("J82").Value = 28057.143
$ws.Range("K82").Value = 3652.3333
$ws.Range("L82").Value = 28057.143
$ws.Range("M82").Value = -3269.3333
$ws.Range("N82").Value = -28823.143
$ws.Range("H85").Value = 23750.412
$ws.Range("I85").Value = 3652.3333
$ws.Range("J85").Value = 28057.143
$ws.Range("K85").Value = 3652.3333
$ws.Range("L85").Value = 28057.143
$ws.Range("M85").Value = -2326.3333
$ws.Range("N85").Value = -30709.143
$ws.Range("H99").Value = 142858350
$ws.Range("I99").Value = 500000100
$ws.Range("K99").Value = 500000100
$ws.Range("M99").Value = -499998602
$ws.Range("H134").Value = 4823.857
$ws.Range("I134").Value = 5680.84
$ws.Range("K134").Value = 17042.52
$ws.Range("M134").Value = -14507.52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 21949.334
$ws.Range("J59").Value = 21949.334
$ws.Range("L59").Value = 21949.334
$ws.Range("N59").Value = -24239.334
$ws.Range("H60").Value = 23996.666
$ws.Range("J60").Value = 23996.666
$ws.Range("L60").Value = 23996.666
$ws.Range("N60").Value = -25018.666
$ws.Range("H63").Value = 36650
$ws.Range("I63").Value = 20000
$ws.Range("K63").Value = 20000
$ws.Range("M63").Value = -19314
$ws.Range("H66").Value = 36650
$ws.Range("I66").Value = 20000
$ws.Range("K66").Value = 60000
$ws.Range("M66").Value = -56568
$ws.Range("H68").Value = 27575
$ws.Range("J68").Value = 27575
$ws.Range("L68").Value = 27575
$ws.Range("N68").Value = -29073
$ws.Range("H71").Value = 27575
$ws.Range("J71").Value = 27575
$ws.Range("L71").Value = 82725
$ws.Range("N71").Value = -90213
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 334277.66
$ws.Range("I5").Value = 636.1818
$ws.Range("J5").Value = 858571.4399999999
$ws.Range("K5").Value = 1908.5454
$ws.Range("L5").Value = 2575714.32
$ws.Range("M5").Value = -1796.5454
$ws.Range("N5").Value = -2575938.32
$ws.Range("H25").Value = 3725
$ws.Range("J25").Value = 4633.3335
$ws.Range("L25").Value = 13900.0005
$ws.Range("N25").Value = -14238.0005
$ws.Range("H30").Value = 3725
$ws.Range("J30").Value = 4633.3335
$ws.Range("L30").Value = 13900.0005
$ws.Range("N30").Value = -14104.0005
$ws.Range("H135").Value = 334277.66
$ws.Range("I135").Value = 636.1818
$ws.Range("J135").Value = 858571.4399999999
$ws.Range("K135").Value = 5725.6362
$ws.Range("L135").Value = 7727142.959999999
$ws.Range("M135").Value = -3190.6362
$ws.Range("N135").Value = -7732212.959999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10000000
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H64").Value = 32635.5
$ws.Range("J64").Value = 32635.5
$ws.Range("L64").Value = 32635.5
$ws.Range("N64").Value = -33131.5
$ws.Range("H67").Value = 32635.5
$ws.Range("J67").Value = 32635.5
$ws.Range("L67").Value = 32635.5
$ws.Range("N67").Value = -34351.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 45005
$ws.Range("J5").Value = 45005
$ws.Range("L5").Value = 45005
$ws.Range("N5").Value = -45231
$ws.Range("H64").Value = 40000
$ws.Range("J64").Value = 40000
$ws.Range("L64").Value = 40000
$ws.Range("N64").Value = -40450
$ws.Range("H67").Value = 40000
$ws.Range("J67").Value = 40000
$ws.Range("L67").Value = 40000
$ws.Range("N67").Value = -41560
$ws.Range("H93").Value = 62525500
$ws.Range("I93").Value = 33917
$ws.Range("K93").Value = 33917
$ws.Range("M93").Value = -32669
$ws.Range("H136").Value = 9924
$ws.Range("I136").Value = 12789.728
$ws.Range("J136").Value = 7297.0835
$ws.Range("K136").Value = 38369.18399999999
$ws.Range("L136").Value = 21891.2505
$ws.Range("M136").Value = -35819.18399999999
$ws.Range("N136").Value = -26991.2505
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 7220
$ws.Range("J22").Value = 7220
$ws.Range("L22").Value = 7220
$ws.Range("N22").Value = -7806
$ws.Range("H64").Value = 25071.334
$ws.Range("J64").Value = 25071.334
$ws.Range("L64").Value = 25071.334
$ws.Range("N64").Value = -25567.334
$ws.Range("H67").Value = 25071.334
$ws.Range("J67").Value = 25071.334
$ws.Range("L67").Value = 25071.334
$ws.Range("N67").Value = -26787.334
